{"js": "// ECU09 \"Registrar Devoluci\u00f3n de Pr\u00e9stamo\" \u2014 Flujo Alterno re-ordering:\n//   1) The \"El sistema modifica el atributo Estado de Alumno...\" step moves up,\n//      ahead of the \"mostrar\u00e1 un MSG...\" / \"El EA selecciona el bot\u00f3n OK.\" steps.\n//   2) The closing step now reads \"...y finaliza el caso de uso.\" instead of\n//      \"...y continua en el paso 27.\".\n//   3) The \"El EA esta logeado al sistema.\" pre-condition keeps the same text\n//      (its run-splitting from the spell-check markup is just simplified).\n\nconst MSG_TEXT =\n  \"El sistema mostrar\u00e1 un MSG \\u201CEl alumno ha devuelto el pr\u00e9stamo fuera del plazo m\u00e1ximo de entrega, se le inhabilitar\u00e1 el acceso al sistema por los pr\u00f3ximos (cantidad de d\u00edas que se pas\u00f3 del plazo de entrega) d\u00edas\\u201D adem\u00e1s del bot\u00f3n OK.\";\nconst OK_TEXT = \"El EA selecciona el bot\u00f3n OK.\";\nconst STATE_TEXT =\n  \"El sistema modifica el atributo Estado de Alumno a \\u201CInhabilitado\\u201D por la cantidad de d\u00edas que pas\u00f3 del plazo de entrega.\";\nconst DISABLE_PREFIX =\n  \"El sistema inhabilita el acceso al Alumno al sistema para posibles futuros pr\u00e9stamos por la cantidad de d\u00edas que pas\u00f3 del plazo de entrega y \";\nconst DISABLE_NEW_SUFFIX = \"finaliza el caso de uso.\";\nconst LOGIN_TEXT = \"El EA esta logeado al sistema.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet msgPara = null;\nlet okPara = null;\nlet statePara = null;\nlet disablePara = null;\nlet loginPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (msgPara === null && t.indexOf(\"mostrar\u00e1 un MSG\") !== -1) {\n    msgPara = p;\n  } else if (okPara === null && t === OK_TEXT) {\n    okPara = p;\n  } else if (statePara === null && t.indexOf(\"modifica el atributo Estado de Alumno\") !== -1) {\n    statePara = p;\n  } else if (disablePara === null && t.indexOf(\"inhabilita el acceso al Alumno al sistema\") !== -1) {\n    disablePara = p;\n  } else if (loginPara === null && t.indexOf(\"logeado\") !== -1) {\n    loginPara = p;\n  }\n}\n\nif (!msgPara || !okPara || !statePara || !disablePara) {\n  throw new Error(\"Could not locate all target paragraphs for the F.Alterno re-order.\");\n}\n\n// Re-order the three steps: the paragraph that used to show the MSG now holds\n// the \"modifica el atributo\" text, the one that used to say \"El EA selecciona\n// el bot\u00f3n OK.\" now holds the MSG text, and the one that used to hold the\n// \"modifica el atributo\" text now says \"El EA selecciona el bot\u00f3n OK.\"\nmsgPara.insertText(STATE_TEXT, Word.InsertLocation.replace);\nokPara.insertText(MSG_TEXT, Word.InsertLocation.replace);\nstatePara.insertText(OK_TEXT, Word.InsertLocation.replace);\n\n// Update the closing sentence of the final alternate-flow step.\ndisablePara.insertText(DISABLE_PREFIX + DISABLE_NEW_SUFFIX, Word.InsertLocation.replace);\n\nif (loginPara) {\n  loginPara.insertText(LOGIN_TEXT, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# ECU09 \"Registrar Devoluci\u00f3n de Pr\u00e9stamo\" \u2014 Flujo Alterno re-ordering:\n#   1) The \"El sistema modifica el atributo Estado de Alumno...\" step moves up,\n#      ahead of the \"mostrar\u00e1 un MSG...\" / \"El EA selecciona el bot\u00f3n OK.\" steps.\n#   2) The closing step now reads \"...y finaliza el caso de uso.\" instead of\n#      \"...y continua en el paso 27.\".\n#   3) The \"El EA esta logeado al sistema.\" pre-condition keeps the same text\n#      (its run-splitting from the spell-check markup is just simplified).\n\n$d = $word.ActiveDocument\n\n$msgText = \"El sistema mostrar\u00e1 un MSG \u201cEl alumno ha devuelto el pr\u00e9stamo fuera del plazo m\u00e1ximo de entrega, se le inhabilitar\u00e1 el acceso al sistema por los pr\u00f3ximos (cantidad de d\u00edas que se pas\u00f3 del plazo de entrega) d\u00edas\u201d adem\u00e1s del bot\u00f3n OK.\"\n$okText = \"El EA selecciona el bot\u00f3n OK.\"\n$stateText = \"El sistema modifica el atributo Estado de Alumno a \u201cInhabilitado\u201d por la cantidad de d\u00edas que pas\u00f3 del plazo de entrega.\"\n$disablePrefix = \"El sistema inhabilita el acceso al Alumno al sistema para posibles futuros pr\u00e9stamos por la cantidad de d\u00edas que pas\u00f3 del plazo de entrega y \"\n$disableNewSuffix = \"finaliza el caso de uso.\"\n$loginText = \"El EA esta logeado al sistema.\"\n\n$msgPara = $null\n$okPara = $null\n$statePara = $null\n$disablePara = $null\n$loginPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($null -eq $msgPara -and $t -like \"*mostrar\u00e1 un MSG*\") {\n        $msgPara = $p\n    } elseif ($null -eq $okPara -and $t -eq $okText) {\n        $okPara = $p\n    } elseif ($null -eq $statePara -and $t -like \"*modifica el atributo Estado de Alumno*\") {\n        $statePara = $p\n    } elseif ($null -eq $disablePara -and $t -like \"*inhabilita el acceso al Alumno al sistema*\") {\n        $disablePara = $p\n    } elseif ($null -eq $loginPara -and $t -like \"*logeado*\") {\n        $loginPara = $p\n    }\n}\n\nif ($null -eq $msgPara -or $null -eq $okPara -or $null -eq $statePara -or $null -eq $disablePara) {\n    throw \"Could not locate all target paragraphs for the F.Alterno re-order.\"\n}\n\n# Replace each whole paragraph's text. Assigning straight to\n# `$paragraph.Range.Text` only overwrites the paragraph's first run (the\n# stale trailing runs survive), so re-materialize the range from explicit\n# Start/End offsets before setting .Text \u2014 that replaces the full span.\n# (Paragraph objects stay live/re-anchored as earlier edits shift offsets,\n# so Start/End are re-read fresh for every assignment below.)\n\n# Re-order the three steps: the paragraph that used to show the MSG now holds\n# the \"modifica el atributo\" text, the one that used to say \"El EA selecciona\n# el bot\u00f3n OK.\" now holds the MSG text, and the one that used to hold the\n# \"modifica el atributo\" text now says \"El EA selecciona el bot\u00f3n OK.\"\n$r1 = $d.Range($msgPara.Range.Start, $msgPara.Range.End)\n$r1.Text = $stateText\n\n$r2 = $d.Range($okPara.Range.Start, $okPara.Range.End)\n$r2.Text = $msgText\n\n$r3 = $d.Range($statePara.Range.Start, $statePara.Range.End)\n$r3.Text = $okText\n\n# Update the closing sentence of the final alternate-flow step.\n$r4 = $d.Range($disablePara.Range.Start, $disablePara.Range.End)\n$r4.Text = $disablePrefix + $disableNewSuffix\n\nif ($null -ne $loginPara) {\n    $r5 = $d.Range($loginPara.Range.Start, $loginPara.Range.End)\n    $r5.Text = $loginText\n}\n"}
